$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: "(Format is 1..." paragraph - drop the _GoBack bookmark that sat
# between the line break and the "(Format is 1" text, and fold that break
# into the same run as the text that follows it.
# ---------------------------------------------------------------------------

$rngFormat = $d.Content
$rngFormat.Find.Execute("(Format is 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngFormat.Text = [char]11 + "(Format is 1"

$d.Bookmarks("_GoBack").Delete()

$rngOldBreak = $d.Content
$rngOldBreak.Find.Execute("Username (generated automatically) ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngOldBreak.Collapse(0)
$rngOldBreak.MoveEnd(1, 1)
$rngOldBreak.Text = ""

# Re-split "Username (generated automatically) " back into its own run
# (it got folded into the following run by the text edit above).
$splitRng = $d.Content
$splitRng.Find.Execute("Username (generated automatically) ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitRng.Collapse(0)
$d.Bookmarks.Add("TempSplitA", $splitRng)
$d.Bookmarks("TempSplitA").Delete()

# ---------------------------------------------------------------------------
# Part 2: "User's Role" paragraph - add the new criteria text on its own
# line, with the _GoBack bookmark now sitting inside it.
# ---------------------------------------------------------------------------

$rngRole = $d.Content
$rngRole.Find.Execute("User's Role ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngRole.Collapse(0)
$rngRole.Text = [char]11 + "(Create 2 roles, an admin and a standard user. By default a user is a standard user)."

# Split "User's Role " into its own run again.
$splitRng2 = $d.Content
$splitRng2.Find.Execute("User's Role ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitRng2.Collapse(0)
$d.Bookmarks.Add("TempSplitB", $splitRng2)
$d.Bookmarks("TempSplitB").Delete()

# Re-create the _GoBack bookmark at its new home, between "...standard u"
# and "ser. By default...".
$rngGoBack = $d.Content
$rngGoBack.Find.Execute("a standard u", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngGoBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngGoBack)

Write-Output "done"
